$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5220
$ws.Range("K3").Value = 8182
$ws.Range("L3").Value = 5629
$ws.Range("J4").Value = 1878
$ws.Range("L4").Value = 1370
$ws.Range("L5").Value = 337
$ws.Range("K6").Value = 9116
$ws.Range("L6").Value = 4683
$ws.Range("J7").Value = 29355
$ws.Range("L7").Value = 17239

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 130
$ws.Range("L7").Value = 558
$ws.Range("L8").Value = 1141
$ws.Range("L11").Value = 281
$ws.Range("L14").Value = 89
$ws.Range("L15").Value = 132
$ws.Range("L16").Value = 35
$ws.Range("L19").Value = 463
$ws.Range("L20").Value = 426
$ws.Range("L24").Value = 44
$ws.Range("L29").Value = 954
$ws.Range("L31").Value = 170
$ws.Range("L33").Value = 801
$ws.Range("L36").Value = 219
$ws.Range("L38").Value = 19
$ws.Range("L42").Value = 564
$ws.Range("J44").Value = 229
$ws.Range("L48").Value = 222
$ws.Range("L51").Value = 215
$ws.Range("L54").Value = 373
$ws.Range("L61").Value = 19
$ws.Range("L63").Value = 50
$ws.Range("L65").Value = 334
$ws.Range("L66").Value = 48
$ws.Range("L71").Value = 48
$ws.Range("L75").Value = 63
$ws.Range("L79").Value = 469
$ws.Range("L83").Value = 378
$ws.Range("L85").Value = 871
$ws.Range("L89").Value = 249
$ws.Range("L90").Value = 176
$ws.Range("L92").Value = 53
$ws.Range("L96").Value = 196
$ws.Range("L97").Value = 143
$ws.Range("J101").Value = 29355
$ws.Range("L101").Value = 17239

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L3").Value = 20
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 59
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 188
$ws.Range("K3").Value = 258
$ws.Range("L3").Value = 188
$ws.Range("K6").Value = 231
$ws.Range("L6").Value = 131
$ws.Range("L7").Value = 558

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 107
$ws.Range("L3").Value = 86
$ws.Range("L7").Value = 281

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 74
$ws.Range("L7").Value = 249

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 263
$ws.Range("L3").Value = 356
$ws.Range("L6").Value = 181
$ws.Range("L7").Value = 871

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L4").Value = 85
$ws.Range("L6").Value = 293
$ws.Range("L7").Value = 1141

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 119
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 378

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 222
$ws.Range("L3").Value = 279
$ws.Range("L4").Value = 49
$ws.Range("L7").Value = 801

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 106
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 334

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 125
$ws.Range("L6").Value = 66

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 170

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 68
$ws.Range("L3").Value = 91
$ws.Range("L6").Value = 181
$ws.Range("L7").Value = 373

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 363
$ws.Range("L6").Value = 235
$ws.Range("L7").Value = 954

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 164
$ws.Range("L7").Value = 463

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 156
$ws.Range("L7").Value = 564

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 150
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 469

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 134
$ws.Range("L7").Value = 426

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 64
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L5").Value = 8
$ws.Range("L6").Value = 19
